$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('minimize', 'мінімізувати', 'c1')
    ,@('grand', 'грандіозний', 'c1')
    ,@('undertake', 'здійснювати', 'c1')
    ,@('instance', 'екземпляр', 'c1')
    ,@('insufficient', 'недостатній', 'c1')
    ,@('establishment', 'створення', 'c1')
    ,@('hop', 'скакати', 'c1')
    ,@('ally', 'союзник', 'c1')
    ,@('bronze', 'бронза', 'c1')
    ,@('pose', 'позу', 'c1')
    ,@('objectively', 'об''єктивно', 'c1')
    ,@('controversy', 'суперечка', 'c1')
    ,@('reinforce', 'підсилювати', 'c1')
    ,@('ample', 'достатньо', 'c1')
    ,@('inferior', 'нижній', 'c1')
    ,@('overestimate', 'переоцінювати', 'c1')
    ,@('mass', 'маса', 'c1')
    ,@('amateur', 'аматор', 'c1')
    ,@('obligatory', 'обов''язковий', 'c1')
    ,@('initiative', 'ініціатива', 'c1')
    ,@('domination', 'домінування', 'c1')
    ,@('self', 'я', 'c1')
    ,@('toll', 'плата', 'c1')
    ,@('oneself', 'себе', 'c1')
    ,@('empire', 'імперія', 'c1')
    ,@('weaken', 'послабити', 'c1')
    ,@('stable', 'конюшня', 'c1')
    ,@('random', 'випадковий', 'c1')
    ,@('simplicity', 'простота', 'c1')
    ,@('worthy', 'гідний', 'c1')
    ,@('exclude', 'виключати', 'c1')
    ,@('migration', 'міграція', 'c1')
    ,@('subsequent', 'наступний', 'c1')
    ,@('rebellion', 'заколот', 'c1')
    ,@('numerous', 'численний', 'c1')
    ,@('restrict', 'обмежувати', 'c1')
    ,@('sufficiently', 'достатньо', 'c1')
    ,@('somewhat', 'дещо', 'c1')
    ,@('determine', 'визначати', 'c1')
    ,@('deficit', 'дефіцит', 'c1')
    ,@('engage', 'займатися', 'c1')
    ,@('restrain', 'обмежувати', 'c1')
    ,@('subsequently', 'згодом', 'c1')
    ,@('constitute', 'складати', 'c1')
    ,@('sacrifice', 'жертвувати', 'c1')
    ,@('purely', 'суто', 'c1')
    ,@('technically', 'технічно', 'c1')
    ,@('unite', 'об''єднувати', 'c1')
    ,@('inclusion', 'включення', 'c1')
    ,@('mutual', 'взаємний', 'c1')
    ,@('tactic', 'тактика', 'c1')
    ,@('seemingly', 'начебто', 'c1')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

Write-Host "Updated rows 2-53"
